$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current total row (row 8), pushing the
# existing total row (8) and footer row (9) down to 10 and 11.
$ws.Rows("8:9").Insert()

# The newly inserted rows 8 and 9 need the same per-column formatting as
# row 7 (the existing item row). Copy row 7's formatting/values down into
# them; the values copied in will be overwritten immediately below.
$ws.Range("A7:Q7").Copy($ws.Range("A8:Q8"))
$ws.Range("A7:Q7").Copy($ws.Range("A9:Q9"))

# Row 7: first item - BRUFEN 400MG 30 TAB
$ws.Range("A7").Value2 = 1
$ws.Range("C7").Value2 = "BRUFEN 400MG 30 TAB"
$ws.Range("H7").Value2 = "0:2"
$ws.Range("L7").Value2 = "1"
$ws.Range("N7").Value2 = "78.00"
$ws.Range("P7").Value2 = "25.7400"
$ws.Range("Q7").Value2 = "0:1"

# Row 8: second item - EUTHYROX 100MCG 50 TAB.
$ws.Range("A8").Value2 = 2
$ws.Range("C8").Value2 = "EUTHYROX 100MCG 50 TAB."
$ws.Range("H8").Value2 = "2:1"
$ws.Range("L8").Value2 = "1"
$ws.Range("N8").Value2 = "70.00"
$ws.Range("P8").Value2 = "70.0000"
$ws.Range("Q8").Value2 = "1:0"

# Row 9: third item - SILVIRBURN CREAM 250 GM (previously row 7's data)
$ws.Range("A9").Value2 = 3
$ws.Range("C9").Value2 = "SILVIRBURN CREAM 250 GM"
$ws.Range("H9").Value2 = "0:0"
$ws.Range("L9").Value2 = "1"
$ws.Range("N9").Value2 = "38.00"
$ws.Range("P9").Value2 = "38.0000"
$ws.Range("Q9").Value2 = "1:0"

# Row 10 (previously row 8): update the running total of selling prices.
$ws.Range("P10").Value2 = 133.74000000000001

# Row 11 (previously row 9): update the generated timestamp.
$ws.Range("A11").Value2 = "Friday, 25 July, 2025 3:34 PM"
